# Add a new prepareStatement test case (ps_dml_004) that exercises
# inserting multiple records in a single prepared statement, and
# relabel the existing single-record insert case (ps_dml_003) so the
# two are distinguishable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 4 (ps_dml_003): clarify title as "single record insert" ----
$ws.Range("C4").Value = '插入单条语句中使用prepareStatement'

# ---- Row 5 (new): ps_dml_004, "multiple record insert" ----
$ws.Range("A5").Value = "ps_dml_004"
$ws.Range("B5").Value = "y"
$ws.Range("C5").Value = '插入多条语句中使用prepareStatement'
$ws.Range("D5").Value = "prepareStatement"
$ws.Range("F5").Value = "schema17"
$ws.Range("H5").Value = 'insert into $schema17 values(?,?,?,?,?,?,?,?,?,?,?,?),(?,?,?,?,?,?,?,?,?,?,?,?),(?,?,?,?,?,?,?,?,?,?,?,?)'
$ws.Range("I5").Value = "1001,JDK_Home,33,35890926187,456.99,98472345827.1299,Tust Plaza 01,19891203,162530,20230320183000,100812,false,1002,zhangsan,-18,88,2.5,12.3,shanghai,20150910,034510,20011111180507,110586,true,1003,awJDs,1,-127,1000.0,-1453.9999,pingYang1,19611001,190000,20101001020202,210092,true"
$ws.Range("J5").Value = "Integer,Varchar,Integer,Bigint,Float,Double,Varchar,Date,Time,Timestamp,Varchar,Boolean,Integer,Varchar,Integer,Bigint,Float,Double,Varchar,Date,Time,Timestamp,Varchar,Boolean,Integer,Varchar,Integer,Bigint,Float,Double,Varchar,Date,Time,Timestamp,Varchar,Boolean"
$ws.Range("K5").Value = "3"
$ws.Range("L5").Value = 'select * from $schema17'
$ws.Range("M5").Value = "src/test/resources/io.dingodb.test/testdata/cases/prepareStatement/expectedresult/ps_dml_004.csv"
$ws.Range("N5").Value = "csv_containsAll"

# Column M ("Expected_result") uses a fill alignment, matching the
# other data rows in the sheet.
$ws.Range("M5").HorizontalAlignment = 5

# Update the view: the previously selected cell (I4) is no longer the
# point of interest; move the selection as in the authored change.
$ws.Range("F24").Select()
